$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.355.77'
$ws.Range("E2").Value = '  +4.36%  '
$ws.Range("D3").Value = '2.046.80'
$ws.Range("E3").Value = '  +3.01%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.31'
$ws.Range("E5").Value = '  +3.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.653'
$ws.Range("E6").Value = '  +2.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '65.46'
$ws.Range("E7").Value = '  +9.80%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.398'
$ws.Range("E9").Value = '  +8.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '59.68'
$ws.Range("E10").Value = '  +0.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0779'
$ws.Range("E11").Value = '  +5.02%  '
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.926'
$ws.Range("E13").Value = '  -2.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.91'
$ws.Range("E14").Value = '  +26.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.80'
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").Value = '2.350.47'
$ws.Range("E16").Value = '  +3.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.69'
$ws.Range("E17").Value = '  +6.57%  '
$ws.Range("D18").Value = '2.058.54'
$ws.Range("E18").Value = '  +3.66%  '
$ws.Range("D19").Value = '37.243.56'
$ws.Range("E19").Value = '  +4.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '73.61'
$ws.Range("E20").Value = '  +2.51%  '
$ws.Range("D21").Value = '0.0₃0884'
$ws.Range("E21").Value = '  +3.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.47'
$ws.Range("E22").Value = '  +4.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.56'
$ws.Range("E23").Value = '  +2.47%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.62'
$ws.Range("E25").Value = '  +1.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.38'
$ws.Range("E26").Value = '  +4.54%  '
$ws.Range("E27").Value = '  +8.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.75'
$ws.Range("E28").Value = '  -2.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.99'
$ws.Range("E29").Value = '  +3.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.131'
$ws.Range("E30").Value = '  +34.73%  '
$ws.Range("E31").Value = '  +2.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.18'
$ws.Range("E32").Value = '  +5.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.19'
$ws.Range("E33").Value = '  +5.44%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.70'
$ws.Range("E34").Value = '  +7.06%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0628'
$ws.Range("E35").Value = '  +4.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.39'
$ws.Range("E36").Value = '  -3.67%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.35'
$ws.Range("E37").Value = '  +11.51%  '
$ws.Range("E38").Value = '  +0.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.84'
$ws.Range("E39").Value = '  +2.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.07'
$ws.Range("E40").Value = '  +31.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.103'
$ws.Range("E41").Value = '  +9.56%  '
$ws.Range("E42").Value = '  +4.28%  '
$ws.Range("E43").Value = '  +7.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.73'
$ws.Range("E44").Value = '  +7.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.16'
$ws.Range("E45").Value = '  +5.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0219'
$ws.Range("E46").Value = '  +2.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '96.20'
$ws.Range("E47").Value = '  +2.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.86'
$ws.Range("E48").Value = '  +1.02%  '
$ws.Range("D49").Value = '1.403.67'
$ws.Range("E49").Value = '  +2.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.93'
$ws.Range("E50").Value = '  +1.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '47.70'
$ws.Range("E51").Value = '  +1.04%  '
